$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to show one record per row with columns:
#   Nombre Completo | Fecha Nacimiento | Nacionalidad | Cedula
# It is being reworked to expose the new "more details" fields returned by
# the details page (requested with the same session), so the header row now
# has 7 columns instead of 4.

# New header row values (A1:G1)
$ws.Range("A1").Value = "Cedula"
$ws.Range("B1").Value = "Nombre:"
$ws.Range("C1").Value = "Primer Apellido:"
$ws.Range("D1").Value = "Segundo Apellido:"
$ws.Range("E1").Value = "Fecha de Nacimiento:"
$ws.Range("F1").Value = "Nacionalidad:"
$ws.Range("G1").Value = "Fallecido/a:"

# Carry the existing bold/bordered/centered header style (style used by
# A1:D1) over to the 3 newly added header cells E1:G1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1:G1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Resize all 7 columns to fit their new header text, same as the original
# columns were already best-fit sized.
$ws.Range("A1:G1").EntireColumn.AutoFit() | Out-Null

# Reflect the new, larger data extent (rows 2-38) in the sheet's selection,
# matching what is shown after the page was re-scraped with more rows.
$ws.Range("A2:G38").Select() | Out-Null
